$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-14 Tuesday" "2025-01-15 Wednesday"

Replace-Text "39×85=3315" "34×39=1326"
Replace-Text "92×74=6808" "62×66=4092"
Replace-Text "97×47=4559" "51×61=3111"
Replace-Text "78×21=1638" "63×41=2583"
Replace-Text "84×30=2520" "76×32=2432"

Replace-Text "95×97=9215" "67×26=1742"
Replace-Text "57×19=1083" "63×52=3276"
Replace-Text "71×55=3905" "28×56=1568"
Replace-Text "62×71=4402" "37×44=1628"
Replace-Text "99×54=5346" "60×23=1380"

Replace-Text "14×50=700" "87×58=5046"
Replace-Text "41×47=1927" "67×59=3953"
Replace-Text "86×76=6536" "86×38=3268"
Replace-Text "82×11=902" "90×53=4770"
Replace-Text "56×63=3528" "48×77=3696"

Replace-Text "45×56=2520" "72×78=5616"
Replace-Text "23×13=299" "43×22=946"
Replace-Text "99×88=8712" "45×99=4455"
Replace-Text "93×54=5022" "44×52=2288"
Replace-Text "31×33=1023" "59×48=2832"

Replace-Text "11×64=704" "33×24=792"
Replace-Text "35×93=3255" "50×97=4850"
Replace-Text "56×34=1904" "81×75=6075"
Replace-Text "45×87=3915" "54×67=3618"
Replace-Text "51×48=2448" "79×75=5925"
